$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Build the "Heading1" / "Heading1Char" style pair (mirrors the
#    built-in Word "Heading 1" style plus its linked character style)
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Style = $d.Styles(-2)          # wdStyleHeading1 -> creates styleId "Heading1"
$hs = $d.Styles(-2)
$hs.NameLocal = "heading 1"
$hs.NextParagraphStyle = "Normal"
$hs.LinkStyle = "Heading1Char"
$hs.ParagraphFormat.SpaceBefore = 12
$hs.Font.Size = 16
$hs.Font.SizeBi = 16
$hs.Font.Color = -738148353              # theme accent1, 25% darker (matches themeShade "BF")

$hc = $d.Styles.Add("Heading 1 Char", 2) # wdStyleTypeCharacter
$hc.BaseStyle = "DefaultParagraphFont"
$hc.Priority = 9
$hc.LinkStyle = "Heading1"
$hc.Font.Size = 16
$hc.Font.SizeBi = 16
$hc.Font.Color = -738148353

# ------------------------------------------------------------------
# 2. Replace the whole document body with the new abstract text
# ------------------------------------------------------------------
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:spacing w:before="0" w:after="60"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
        <w:t>Clear Creek</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve">The U.S. Fish and Wildlife Service continued its juvenile salmonid monitoring project in Clear Creek, Shasta County, California, that began in December 1998 using a rotary screw trap (RST) located at river mile (RM) 1.7. This monitoring project has three primary objectives: (1) calculate annual juvenile passage indices for Chinook Salmon </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:i/>
        </w:rPr>
        <w:t>Oncorhynchus tshawytscha</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> and Rainbow Trout/steelhead </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:i/>
        </w:rPr>
        <w:t>O. mykiss</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> (we use the term </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>O. mykiss</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> to refer to both the stream resident (Rainbow Trout) and anadromous (steelhead) life histories because of the difficulties in differentiating the anadromous and resident forms in the field) for inter-year comparisons and analyses of the effectiveness of stream restoration activities; (2) obtain juvenile salmonid life history information including size, timing, and potential factors limiting survival; and (3) collect biological samples from juvenile salmonids. Chinook Salmon run classifications show that all Chinook Salmon runs known to the Sacramento River were captured; however, due to overlapping spawning times of spring-run and fall-run Chinook Salmon it was problematic to index the juvenile passage using only the RST at RM 1.7. Since 2003 a separation weir has been used to isolate adult spring-run Chinook Salmon from adult fall-run Chinook Salmon; therefore, a second RST was added above the weir at RM 8.4.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:spacing w:before="120" w:after="60"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
        <w:t>Battle Creek</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="240"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve">The U.S. Fish and Wildlife Service continued its juvenile salmonid monitoring project in Battle Creek, Shasta and Tehama Counties, California, that began in December 1998 using </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> RST located at river mile (RM) 6.2. The upper Battle Creek RST objectives are to (1) calculate annual juvenile passage indices for spring-run, late-fall run, and winter-run Chinook Salmon; Rainbow Trout/steelhead </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:i/>
        </w:rPr>
        <w:t>O. mykiss</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> (we use the term </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>O. mykiss</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> to refer to both the stream resident (Rainbow Trout) and anadromous (steelhead) life histories because of the difficulties in differentiating the anadromous and resident forms in the field)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t>(2) obtain juvenile salmonid life history information including size, timing, and potential factors limiting survival; and (3) collect biological samples from juvenile salmonids. Chinook Salmon run designations in the Sacramento River watershed were developed using length-at-date criteria for Sacramento River fall-run Chinook Salmon. These designations may not be accurate when used for salmon captured in tributaries of the Sacramento River; these streams have variable and overlapping run timing of spring-run and fall-run Chinook Salmon. After 2008 all fall-run designated Chinook Salmon were reassigned as spring-run in the spring-run Chinook Salmon passage indices as preventative measures were taken to reduce the potential for fall-run Chinook Salmon to enter the upper watershed (above the Coleman National Fish Hatchery barrier weir).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:spacing w:before="120" w:after="60"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="auto"/>
        </w:rPr>
        <w:t>Mark–recapture Trials</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="240"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="SimSun" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve">Since the RST only captures fish from a small portion of the creek cross section, it is necessary to implement a method to project the RST catch numbers to portions of the creek outside of the RST capture zone. Mark–recapture trials were attempted, when creek conditions were conducive for such studies and when fish numbers allowed, to determine the efficiency of the RSTs to catch juvenile salmonids moving downstream during a given </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="SimSun" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t>time period</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve">. If one or no mark–recaptures trials were conducted, then the historic monthly </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>and global averages were used to calculate that passage indices. Mark-recapture–trials on Battle Creek are conducted using hatchery-origin fish (from Coleman National Fish Hatchery), while on Clear Creek natural-origin fish are used. Fish used for the trials were either single or dual marked. For the single mark trials the fish were stained with Bismarck brown-Y. Dual marked fish clipped with surgical scalpels, removing an area of approximately 1–2 mm</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t xml:space="preserve"> from the tip of either the upper or lower lobe of the caudal fin and then the fish were stained as above.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
      </w:pPr>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
